$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep Text storage like the source data,
# since plain numeric-looking strings would otherwise be auto-coerced to numbers.
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "29.587.82"
$ws.Cells.Item(2, 5).Value = "  +0.26%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.923.11"
$ws.Cells.Item(3, 5).Value = "  +0.47%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.013"
$ws.Cells.Item(4, 5).Value = "  +0.56%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "325.68"
$ws.Cells.Item(5, 5).Value = "  -0.12%  "
$ws.Cells.Item(6, 5).Value = "  +0.39%  "
$ws.Cells.Item(8, 5).Value = "  -0.31%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.08222"
$ws.Cells.Item(9, 5).Value = "  +1.06%  "
$ws.Cells.Item(10, 5).Value = "  -0.17%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "23.68"
$ws.Cells.Item(11, 5).Value = "  +1.27%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "1.930.70"
$ws.Cells.Item(12, 5).Value = "  +0.65%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "6.073"
$ws.Cells.Item(13, 5).Value = "  +1.23%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "7.263"
$ws.Cells.Item(14, 5).Value = "  +1.79%  "
$ws.Cells.Item(15, 5).Value = "  +1.55%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.06861"
$ws.Cells.Item(16, 5).Value = "  +1.22%  "
$ws.Cells.Item(17, 5).Value = "  +0.34%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.00001036"
$ws.Cells.Item(18, 5).Value = "  -0.28%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "17.63"
$ws.Cells.Item(19, 5).Value = "  -0.37%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "1.011"
$ws.Cells.Item(20, 5).Value = "  +0.42%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "29.587.62"
$ws.Cells.Item(21, 5).Value = "  +0.22%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "5.685"
$ws.Cells.Item(22, 5).Value = "  +1.06%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "11.94"
$ws.Cells.Item(23, 5).Value = "  +1.30%  "
$ws.Cells.Item(24, 5).Value = "  +0.01%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.168.94"
$ws.Cells.Item(25, 5).Value = "  +0.93%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "155.83"
$ws.Cells.Item(26, 5).Value = "  +0.28%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "6.452"
$ws.Cells.Item(27, 5).Value = "  +0.32%  "
$ws.Cells.Item(28, 5).Value = "  -0.19%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "2.091"
$ws.Cells.Item(29, 5).Value = "  -0.27%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "120.73"
$ws.Cells.Item(30, 5).Value = "  +0.85%  "
$ws.Cells.Item(31, 5).Value = "  -1.62%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.09628"
$ws.Cells.Item(32, 5).Value = "  +0.71%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "5.619"
$ws.Cells.Item(33, 5).Value = "  +2.04%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "3.563"
$ws.Cells.Item(34, 5).Value = "  -0.11%  "
$ws.Cells.Item(35, 5).Value = "  -1.05%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.06352"
$ws.Cells.Item(36, 5).Value = "  +4.14%  "
$ws.Cells.Item(37, 5).Value = "  +0.99%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "1.184"
$ws.Cells.Item(38, 5).Value = "  +0.76%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.5957"
$ws.Cells.Item(39, 5).Value = "  +0.34%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "10.75"
$ws.Cells.Item(40, 5).Value = "  -0.09%  "
$ws.Cells.Item(41, 2).Value = "FraxShare"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "7.876"
$ws.Cells.Item(41, 5).Value = "  -1.24%  "
$ws.Cells.Item(42, 2).Value = "Algorand"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.1851"
$ws.Cells.Item(42, 5).Value = "  -0.40%  "
$ws.Cells.Item(43, 2).Value = "RenderToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "2.453"
$ws.Cells.Item(43, 5).Value = "  -0.37%  "
$ws.Cells.Item(44, 2).Value = "WEMIXToken"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "1.286"
$ws.Cells.Item(44, 5).Value = "  +3.22%  "
$ws.Cells.Item(45, 2).Value = "EnergySwap"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "12.38"
$ws.Cells.Item(45, 5).Value = "  -0.21%  "
$ws.Cells.Item(46, 2).Value = "Cronos"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.07502"
$ws.Cells.Item(46, 5).Value = "  -2.49%  "
$ws.Cells.Item(47, 2).Value = "Decentraland"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.5566"
$ws.Cells.Item(47, 5).Value = "  -0.18%  "
$ws.Cells.Item(48, 2).Value = "NEARProtocol"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "1.951"
$ws.Cells.Item(48, 5).Value = "  +0.35%  "
$ws.Cells.Item(49, 2).Value = "Quant"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "119.55"
$ws.Cells.Item(49, 5).Value = "  +3.40%  "
$ws.Cells.Item(50, 2).Value = "MXToken"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "2.435"
$ws.Cells.Item(50, 5).Value = "  +3.38%  "
$ws.Cells.Item(51, 2).Value = "Aave"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "72.13"
$ws.Cells.Item(51, 5).Value = "  -0.83%  "
